$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure these Price cells remain text (not auto-converted to numbers),
# matching the original inline-string cell type.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = '35.524.18'
$ws.Range("E2").Value = '  +1.30%  '

$ws.Range("D3").Value = '1.907.53'
$ws.Range("E3").Value = '  +2.66%  '

$ws.Range("E4").Value = '  +0.30%  '

$ws.Range("D5").Value = '246.76'
$ws.Range("E5").Value = '  +3.53%  '

$ws.Range("D6").Value = '0.657'
$ws.Range("E6").Value = '  +5.61%  '

$ws.Range("E7").Value = '  +0.18%  '

$ws.Range("D8").Value = '42.05'
$ws.Range("E8").Value = '  -0.94%  '

$ws.Range("E9").Value = '  +3.95%  '

$ws.Range("D10").Value = '49.05'
$ws.Range("E10").Value = '  +5.22%  '

$ws.Range("D11").Value = '0.0713'
$ws.Range("E11").Value = '  +2.31%  '

$ws.Range("D12").Value = '0.1000'
$ws.Range("E12").Value = '  +0.86%  '

$ws.Range("D13").Value = '2.182.10'
$ws.Range("E13").Value = '  +2.60%  '

$ws.Range("D14").Value = '12.39'
$ws.Range("E14").Value = '  +8.62%  '

$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").Value = '0.700'
$ws.Range("E15").Value = '  +3.31%  '

$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '1.902.70'
$ws.Range("E16").Value = '  +2.25%  '

$ws.Range("D17").Value = '4.85'
$ws.Range("E17").Value = '  +3.03%  '

$ws.Range("D18").Value = '35.530.79'
$ws.Range("E18").Value = '  +1.40%  '

$ws.Range("D19").Value = '72.23'
$ws.Range("E19").Value = '  +2.76%  '

$ws.Range("D20").Value = '0.0₃0833'
$ws.Range("E20").Value = '  +4.40%  '

$ws.Range("D21").Value = '243.88'
$ws.Range("E21").Value = '  +1.18%  '

$ws.Range("D22").Value = '12.64'
$ws.Range("E22").Value = '  +4.03%  '

$ws.Range("D23").Value = '4.84'
$ws.Range("E23").Value = '  +2.05%  '

$ws.Range("E24").Value = '  +0.17%  '

$ws.Range("E25").Value = '  +1.33%  '

$ws.Range("D26").Value = '2.21'
$ws.Range("E26").Value = '  +13.59%  '

$ws.Range("D27").Value = '171.25'
$ws.Range("E27").Value = '  -0.14%  '

$ws.Range("D28").Value = '8.53'
$ws.Range("E28").Value = '  +7.71%  '

$ws.Range("D29").Value = '18.09'
$ws.Range("E29").Value = '  +2.12%  '

$ws.Range("D30").Value = '0.129'
$ws.Range("E30").Value = '  +3.22%  '

$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").Value = '0.968'
$ws.Range("E31").Value = '  +22.44%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '4.17'
$ws.Range("E32").Value = '  +4.26%  '

$ws.Range("D33").Value = '0.0571'
$ws.Range("E33").Value = '  +2.31%  '

$ws.Range("D34").Value = '4.23'
$ws.Range("E34").Value = '  +5.45%  '

$ws.Range("E35").Value = '  +0.29%  '

$ws.Range("D36").Value = '1.74'
$ws.Range("E36").Value = '  +7.52%  '

$ws.Range("E37").Value = '  +0.23%  '

$ws.Range("E38").Value = '  +2.73%  '

$ws.Range("E39").Value = '  +2.63%  '

$ws.Range("D40").Value = '92.48'
$ws.Range("E40").Value = '  +0.86%  '

$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").Value = '0.0205'
$ws.Range("E41").Value = '  +1.21%  '

$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").Value = '0.0638'
$ws.Range("E42").Value = '  +17.41%  '

$ws.Range("D43").Value = '15.67'
$ws.Range("E43").Value = '  +5.09%  '

$ws.Range("D44").Value = '1.348.83'
$ws.Range("E44").Value = '  -0.40%  '

$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").Value = '2.39'
$ws.Range("E45").Value = '  +1.93%  '

$ws.Range("B46").Value = 'MultiversX'
$ws.Range("C46").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D46").Value = '48.60'
$ws.Range("E46").Value = '  +40.87%  '

$ws.Range("D47").Value = '12.63'
$ws.Range("E47").Value = '  -1.54%  '

$ws.Range("E48").Value = '  +0.07%  '

$ws.Range("E49").Value = '  +0.28%  '

$ws.Range("D50").Value = '6.61'
$ws.Range("E50").Value = '  +3.09%  '

$ws.Range("D51").Value = '2.093.82'
$ws.Range("E51").Value = '  +2.63%  '
